# Update canine keywords on the CypherOutput sheet (data row 2):
#  - Case ID:            NCATS-COP01CCB010072 -> NCATS-COP01-CCB010072
#  - Diagnosis:           Bone sarcomas :: Osteosarcoma (appendicular) -> Osteosarcoma
#  - Stage of Disease:    (empty) -> Unknown

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CypherOutput")

$ws.Range("A2").Value = "NCATS-COP01-CCB010072"
$ws.Range("E2").Value = "Osteosarcoma"
$ws.Range("F2").Value = "Unknown"
